$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$srcStyle = $ws.Range("C10").Style
$ws.Range("B5").Style = $srcStyle
